$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Add new column G header and values (test for parsing NaN)
$ws.Range("G1").Value = "not available"
$ws.Range("G2").Value = "NaN"
$ws.Range("G3").Value = "NaN"

# Set column G width similar to the authored change (closest value the
# COM width-quantization model can reach to the authored 16.734375)
$ws.Columns.Item(7).ColumnWidth = 15.833333333333334

# Update selection to mirror the recorded view state after the edit
$ws.Range("H12").Select() | Out-Null
